$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = -11.16279999999998
$ws.Range("A12").Value = -21.90460000000002
$ws.Range("C12").Value = -13.1366
$ws.Range("C14").Value = -11.87759999999999
$ws.Range("C22").Value = -10.81159999999999
